$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- uInt32Array (mirrors the existing uDoubleArray block at row 12) ---
$ws.Range("A34").Value = "uInt32Array"
$ws.Range("B34:E34").FormulaArray = "=_xll.uInt32Array(B35:C36)"

$ws.Range("B35").Value = 1
$ws.Range("C35").Value = 2
$ws.Range("B36").Value = 3
$ws.Range("C36").Value = 4

# --- uInt32Matrix (mirrors the existing uDateTimeMatrix block at row 15) ---
$ws.Range("A37").Value = "uInt32Matrix"
$ws.Range("B37:C38").FormulaArray = "=_xll.uInt32Matrix(B35:C36)"

# Keep the new block's selection/view consistent with where the author
# finished editing.
$ws.Range("B38").Select()
